$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "tar"
$ws.Range("E2").Value = "N/A"

[void]$ws.Range("E4").Select()
